$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (51) of profit data for 01/14/2026, following the
# pattern of the existing rows in the sheet.
$row = 51

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "01/14/2026"
$ws.Cells.Item($row, 1).ClearFormats()
$ws.Cells.Item($row, 2).Value = 13305.92
$ws.Cells.Item($row, 3).Value = 0.2172232585151117
$ws.Cells.Item($row, 4).Value = 0.7827767414848883
$ws.Cells.Item($row, 5).Value = -107.35
$ws.Cells.Item($row, 6).Value = -16.91
$ws.Cells.Item($row, 7).Value = -20276.2
$ws.Cells.Item($row, 8).Value = -66.06
$ws.Cells.Item($row, 9).Value = -262.5
$ws.Cells.Item($row, 10).Value = -8.33
